# Add data for 2021-11-16 (update "through November 07" -> "through November 08")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update the header label / title cell (shared string for column B header)
$ws.Name = "Through 2021-11-08"
$ws.Range("B1").Value = "November 2021 (through November 08)"

# Update existing cell values (carjacking counts changed for the "through Nov 08" column)
$ws.Range("M2").Value = 5
$ws.Range("BE3").Value = 2
$ws.Range("AT4").Value = 5
$ws.Range("BP4").Value = 2
$ws.Range("B6").Value = 5
$ws.Range("M6").Value = 2
$ws.Range("AT9").Value = 6
$ws.Range("AT44").Value = 2
$ws.Range("M47").Value = 2

# Fill in newly populated (previously empty) cells
$ws.Range("M10").Value = 1
$ws.Range("AI12").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("X21").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("B24").Value = 1
$ws.Range("BE24").Value = 1
$ws.Range("B29").Value = 1
$ws.Range("AT42").Value = 1
$ws.Range("AI66").Value = 1
$ws.Range("BP79").Value = 1
$ws.Range("X80").Value = 1
$ws.Range("AT88").Value = 1
